# "modified job posting dict" -- split the single "Category" column of the
# Job Postings Dict sheet into three columns: Category (B), Field (C) and
# Skills (D, the text that used to live in B).

$wb = $excel.ActiveWorkbook

# --- Row height tweaks on "Manual" and "Professional" sheets (row 1: default -> 13) ---
$wb.Worksheets.Item("Manual").Rows.Item(1).RowHeight = 13
$wb.Worksheets.Item("Professional").Rows.Item(1).RowHeight = 13

# --- "Job Postings Dict" sheet ---
$ws = $wb.Worksheets.Item("Job Postings Dict")

# Capture the existing "skills" text (currently in column B) before overwriting it.
$skills1 = $ws.Range("B1").Value2
$skills2 = $ws.Range("B2").Value2
$skills3 = $ws.Range("B3").Value2
$skills4 = $ws.Range("B4").Value2
$skills5 = $ws.Range("B5").Value2

# Move column B's current formatting over to column D (where the skills text
# will now live) before the B column cells get re-purposed/re-styled.
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Write the skills text into column D.
$ws.Range("D1").Value2 = $skills1
$ws.Range("D2").Value2 = $skills2
$ws.Range("D3").Value2 = $skills3
$ws.Range("D4").Value2 = $skills4
$ws.Range("D5").Value2 = $skills5

# Reset columns B and C to the plain/default style before writing the new
# category + field text into them.
$ws.Range("B1:C5").Style = "Normal"

# Row 1: Data Analyst -> Technical / Data Science
$ws.Range("B1").Value2 = "Technical"
$ws.Range("C1").Value2 = "Data Science"

# Row 2: Construction Worker -> Manual / Construction
$ws.Range("B2").Value2 = "Manual"
$ws.Range("C2").Value2 = "Construction"

# Row 3: Babysitter -> Manual / Child Care
$ws.Range("B3").Value2 = "Manual"
$ws.Range("C3").Value2 = "Child Care"

# Row 4: Doctor -> Professional / Medical
$ws.Range("B4").Value2 = "Professional"
$ws.Range("C4").Value2 = "Medical"

# Row 5: Math Teacher -> Professional / Education
$ws.Range("B5").Value2 = "Professional"
$ws.Range("C5").Value2 = "Education"

# Move the active selection to B7 to match the saved view state
$ws.Range("B7").Select()
